$d = $word.ActiveDocument

# --- Change 1: "Kjetil Raaen" paragraph -> single clean run, no proofErr wrapping ---
$para2 = $d.Paragraphs(2)
$atStart = $para2.Range.Duplicate
$atStart.Collapse(1)
$atStart.InsertParagraphBefore()
$newNameRange = $d.Paragraphs(2).Range.Duplicate
$newNameRange.Collapse(1)
$newNameRange.InsertBefore("Kjetil Raaen")
$oldNamePara = $d.Paragraphs(3)
$oldNamePara.Range.Delete()

# --- Change 2: drop the stray "_GoBack" bookmark after "27" / before "th" ---
$d.Bookmarks("_GoBack").Delete()

# --- Change 3: add explanatory sentence to the ASCII-art paragraph, and
#     re-add the "_GoBack" bookmark there (collapsed at the paragraph end) ---
$asciiPara = $d.Paragraphs(9)
$endOfPara = $asciiPara.Range.Duplicate
$endOfPara.Collapse(0)
$endOfPara.MoveEnd(1, -1)
$endOfPara.InsertAfter("That is, the asci-art is split into 30 by 30 squares of characters, and need to be puzzled together. After each line you will need to add a newline character.")

$bm = $asciiPara.Range.Duplicate
$bm.Collapse(0)
$bm.MoveEnd(1, -1)
$d.Bookmarks.Add("_GoBack", $bm)
